$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values (A1:B32)
$ws.Cells.Item(1, 1).Value = -0.1204449352092638
$ws.Cells.Item(1, 2).Value = 0.12001252009869745
$ws.Cells.Item(2, 1).Value = -0.073905022829290345
$ws.Cells.Item(2, 2).Value = 0.072120932331944942
$ws.Cells.Item(3, 1).Value = 0.058049412124919542
$ws.Cells.Item(3, 2).Value = -0.058541377216563717
$ws.Cells.Item(4, 1).Value = -0.22544931194365248
$ws.Cells.Item(4, 2).Value = 0.22384274668471349
$ws.Cells.Item(5, 1).Value = -0.13974666021111393
$ws.Cells.Item(5, 2).Value = 0.13834849849425357
$ws.Cells.Item(6, 1).Value = -0.084524263644018216
$ws.Cells.Item(6, 2).Value = 0.084417706767952794
$ws.Cells.Item(7, 1).Value = -0.064417707607374197
$ws.Cells.Item(7, 2).Value = 0.064185236026823134
$ws.Cells.Item(8, 1).Value = -0.044185236873445
$ws.Cells.Item(8, 2).Value = 0.04402787814073772
$ws.Cells.Item(9, 1).Value = -0.038027878875413812
$ws.Cells.Item(9, 2).Value = 0.037907760122791068
$ws.Cells.Item(10, 1).Value = -0.031907760864221757
$ws.Cells.Item(10, 2).Value = 0.031897926730650283
$ws.Cells.Item(11, 1).Value = -0.027397927459229976
$ws.Cells.Item(11, 2).Value = 0.027379281696219948
$ws.Cells.Item(12, 1).Value = -0.021379282439069947
$ws.Cells.Item(12, 2).Value = 0.021327536254316826
$ws.Cells.Item(13, 1).Value = -0.015327537002604252
$ws.Cells.Item(13, 2).Value = 0.015318667205735714
$ws.Cells.Item(14, 1).Value = -0.0033186680066759067
$ws.Cells.Item(14, 2).Value = 0.0033183645546497331
$ws.Cells.Item(15, 1).Value = -0.021052275190096559
$ws.Cells.Item(15, 2).Value = 0.021027322158857054
$ws.Cells.Item(16, 1).Value = -0.015027322911285612
$ws.Cells.Item(16, 2).Value = 0.015004094752339547
$ws.Cells.Item(17, 1).Value = -0.0090040955080548102
$ws.Cells.Item(17, 2).Value = 0.0089999992178624311
$ws.Cells.Item(18, 1).Value = -0.03610749798323809
$ws.Cells.Item(18, 2).Value = 0.03609616305802632
$ws.Cells.Item(19, 1).Value = -0.027096163763987491
$ws.Cells.Item(19, 2).Value = 0.027013169852033059
$ws.Cells.Item(20, 1).Value = -0.018013170564117331
$ws.Cells.Item(20, 2).Value = 0.018004192156096721
$ws.Cells.Item(21, 1).Value = -0.0090041928690896
$ws.Cells.Item(21, 2).Value = 0.0089999992863658562
$ws.Cells.Item(22, 1).Value = -0.093938911641354039
$ws.Cells.Item(22, 2).Value = 0.09362908972877193
$ws.Cells.Item(23, 1).Value = -0.084629090446637356
$ws.Cells.Item(23, 2).Value = 0.084125603221800738
$ws.Cells.Item(24, 1).Value = -0.042125604234104053
$ws.Cells.Item(24, 2).Value = 0.041999998982228171
$ws.Cells.Item(25, 1).Value = -0.047824235567151874
$ws.Cells.Item(25, 2).Value = 0.047775855029009051
$ws.Cells.Item(26, 1).Value = -0.041775855746511326
$ws.Cells.Item(26, 2).Value = 0.041718988247655631
$ws.Cells.Item(27, 1).Value = -0.03571898896638892
$ws.Cells.Item(27, 2).Value = 0.035536750389821847
$ws.Cells.Item(28, 1).Value = -0.029536751114974891
$ws.Cells.Item(28, 2).Value = 0.029427225090444153
$ws.Cells.Item(29, 1).Value = -0.017427225870692453
$ws.Cells.Item(29, 2).Value = 0.017393613732483715
$ws.Cells.Item(30, 1).Value = 0.00260638541760283
$ws.Cells.Item(30, 2).Value = -0.0026765887115258735
$ws.Cells.Item(31, 1).Value = 0.017676587904736252
$ws.Cells.Item(31, 2).Value = -0.017714160131898637
$ws.Cells.Item(32, 1).Value = -0.0060005649980512388
$ws.Cells.Item(32, 2).Value = 0.0059999992708945626

# Column B width changed from 14.7109375 to 15.42578125 (character width units)
$ws.Columns.Item(2).ColumnWidth = 14.7
